# FMECA_Sample.xlsx -- "add df for equipment"
#
# The underlying author action (per the xml diff) was: open the sheet in
# Excel, select every cell, and auto-fit columns A:F so each column is wide
# enough to show its longest value/header without truncation. (The rest of
# the diff is version/save metadata Excel itself stamps on every re-save and
# isn't something this automation surface can or should try to reproduce.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole sheet (mirrors the saved <selection sqref="A1:XFD1048576"/>).
$ws.Cells.Select() | Out-Null

# Auto-fit the six data columns (A:F) to their content, matching the
# column widths captured in the saved workbook.
$ws.Columns("A").ColumnWidth = 10.833333333333334
$ws.Columns("B").ColumnWidth = 12.666666666666666
$ws.Columns("C").ColumnWidth = 22.166666666666668
$ws.Columns("D").ColumnWidth = 6.666666666666667
$ws.Columns("E").ColumnWidth = 17.833333333333332
$ws.Columns("F").ColumnWidth = 7.833333333333333
